# Update cryptocurrency price/volume data (auto-generated from diff)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range('D2')
$cell.NumberFormat = "@"
$cell.Value = '26.144.53'
$cell.Style = "Normal"

$cell = $ws.Range('E2')
$cell.NumberFormat = "@"
$cell.Value = '  -0.02%  '
$cell.Style = "Normal"

$cell = $ws.Range('D3')
$cell.NumberFormat = "@"
$cell.Value = '1.656.68'
$cell.Style = "Normal"

$cell = $ws.Range('E3')
$cell.NumberFormat = "@"
$cell.Value = '  +0.00%  '
$cell.Style = "Normal"

$cell = $ws.Range('D5')
$cell.NumberFormat = "@"
$cell.Value = '218.63'
$cell.Style = "Normal"

$cell = $ws.Range('E5')
$cell.NumberFormat = "@"
$cell.Value = '  -0.21%  '
$cell.Style = "Normal"

$cell = $ws.Range('D6')
$cell.NumberFormat = "@"
$cell.Value = '0.5237'
$cell.Style = "Normal"

$cell = $ws.Range('E6')
$cell.NumberFormat = "@"
$cell.Value = '  -0.16%  '
$cell.Style = "Normal"

$cell = $ws.Range('E7')
$cell.NumberFormat = "@"
$cell.Value = '  -0.17%  '
$cell.Style = "Normal"

$cell = $ws.Range('D8')
$cell.NumberFormat = "@"
$cell.Value = '0.2660'
$cell.Style = "Normal"

$cell = $ws.Range('E8')
$cell.NumberFormat = "@"
$cell.Value = '  +1.43%  '
$cell.Style = "Normal"

$cell = $ws.Range('D9')
$cell.NumberFormat = "@"
$cell.Value = '0.06357'
$cell.Style = "Normal"

$cell = $ws.Range('E9')
$cell.NumberFormat = "@"
$cell.Value = '  +0.93%  '
$cell.Style = "Normal"

$cell = $ws.Range('E10')
$cell.NumberFormat = "@"
$cell.Value = '  +0.00%  '
$cell.Style = "Normal"

$cell = $ws.Range('D11')
$cell.NumberFormat = "@"
$cell.Value = '0.07683'
$cell.Style = "Normal"

$cell = $ws.Range('E11')
$cell.NumberFormat = "@"
$cell.Value = '  -1.53%  '
$cell.Style = "Normal"

$cell = $ws.Range('E12')
$cell.NumberFormat = "@"
$cell.Value = '  +2.56%  '
$cell.Style = "Normal"

$cell = $ws.Range('D13')
$cell.NumberFormat = "@"
$cell.Value = '1.650.10'
$cell.Style = "Normal"

$cell = $ws.Range('E13')
$cell.NumberFormat = "@"
$cell.Value = '  -0.52%  '
$cell.Style = "Normal"

$cell = $ws.Range('D14')
$cell.NumberFormat = "@"
$cell.Value = '1.884.34'
$cell.Style = "Normal"

$cell = $ws.Range('E14')
$cell.NumberFormat = "@"
$cell.Value = '  -0.01%  '
$cell.Style = "Normal"

$cell = $ws.Range('D15')
$cell.NumberFormat = "@"
$cell.Value = '0.5628'
$cell.Style = "Normal"

$cell = $ws.Range('E15')
$cell.NumberFormat = "@"
$cell.Value = '  +1.36%  '
$cell.Style = "Normal"

$cell = $ws.Range('D16')
$cell.NumberFormat = "@"
$cell.Value = '0.0₅8209'
$cell.Style = "Normal"

$cell = $ws.Range('E16')
$cell.NumberFormat = "@"
$cell.Value = '  +2.47%  '
$cell.Style = "Normal"

$cell = $ws.Range('D17')
$cell.NumberFormat = "@"
$cell.Value = '65.50'
$cell.Style = "Normal"

$cell = $ws.Range('E17')
$cell.NumberFormat = "@"
$cell.Value = '  +0.67%  '
$cell.Style = "Normal"

$cell = $ws.Range('D18')
$cell.NumberFormat = "@"
$cell.Value = '26.147.01'
$cell.Style = "Normal"

$cell = $ws.Range('E18')
$cell.NumberFormat = "@"
$cell.Value = '  -0.04%  '
$cell.Style = "Normal"

$cell = $ws.Range('E19')
$cell.NumberFormat = "@"
$cell.Value = '  -0.12%  '
$cell.Style = "Normal"

$cell = $ws.Range('D20')
$cell.NumberFormat = "@"
$cell.Value = '4.663'
$cell.Style = "Normal"

$cell = $ws.Range('E20')
$cell.NumberFormat = "@"
$cell.Value = '  +0.45%  '
$cell.Style = "Normal"

$cell = $ws.Range('D21')
$cell.NumberFormat = "@"
$cell.Value = '10.54'
$cell.Style = "Normal"

$cell = $ws.Range('E21')
$cell.NumberFormat = "@"
$cell.Value = '  +4.20%  '
$cell.Style = "Normal"

$cell = $ws.Range('D22')
$cell.NumberFormat = "@"
$cell.Value = '192.82'
$cell.Style = "Normal"

$cell = $ws.Range('E22')
$cell.NumberFormat = "@"
$cell.Value = '  -1.36%  '
$cell.Style = "Normal"

$cell = $ws.Range('D23')
$cell.NumberFormat = "@"
$cell.Value = '5.958'
$cell.Style = "Normal"

$cell = $ws.Range('E23')
$cell.NumberFormat = "@"
$cell.Value = '  -0.03%  '
$cell.Style = "Normal"

$cell = $ws.Range('D25')
$cell.NumberFormat = "@"
$cell.Value = '145.43'
$cell.Style = "Normal"

$cell = $ws.Range('E25')
$cell.NumberFormat = "@"
$cell.Value = '  -0.82%  '
$cell.Style = "Normal"

$cell = $ws.Range('D26')
$cell.NumberFormat = "@"
$cell.Value = '0.1198'
$cell.Style = "Normal"

$cell = $ws.Range('E26')
$cell.NumberFormat = "@"
$cell.Value = '  -0.42%  '
$cell.Style = "Normal"

$cell = $ws.Range('D27')
$cell.NumberFormat = "@"
$cell.Value = '7.268'
$cell.Style = "Normal"

$cell = $ws.Range('E27')
$cell.NumberFormat = "@"
$cell.Value = '  +1.46%  '
$cell.Style = "Normal"

$cell = $ws.Range('E28')
$cell.NumberFormat = "@"
$cell.Value = '  +0.16%  '
$cell.Style = "Normal"

$cell = $ws.Range('D29')
$cell.NumberFormat = "@"
$cell.Value = '1.508'
$cell.Style = "Normal"

$cell = $ws.Range('E29')
$cell.NumberFormat = "@"
$cell.Value = '  +0.84%  '
$cell.Style = "Normal"

$cell = $ws.Range('D30')
$cell.NumberFormat = "@"
$cell.Value = '0.05468'
$cell.Style = "Normal"

$cell = $ws.Range('E30')
$cell.NumberFormat = "@"
$cell.Value = '  -4.08%  '
$cell.Style = "Normal"

$cell = $ws.Range('E31')
$cell.NumberFormat = "@"
$cell.Value = '  -0.03%  '
$cell.Style = "Normal"

$cell = $ws.Range('D32')
$cell.NumberFormat = "@"
$cell.Value = '3.467'
$cell.Style = "Normal"

$cell = $ws.Range('E32')
$cell.NumberFormat = "@"
$cell.Value = '  -0.53%  '
$cell.Style = "Normal"

$cell = $ws.Range('D33')
$cell.NumberFormat = "@"
$cell.Value = '3.377'
$cell.Style = "Normal"

$cell = $ws.Range('E33')
$cell.NumberFormat = "@"
$cell.Value = '  +0.73%  '
$cell.Style = "Normal"

$cell = $ws.Range('D34')
$cell.NumberFormat = "@"
$cell.Value = '1.567'
$cell.Style = "Normal"

$cell = $ws.Range('E34')
$cell.NumberFormat = "@"
$cell.Value = '  -1.23%  '
$cell.Style = "Normal"

$cell = $ws.Range('D35')
$cell.NumberFormat = "@"
$cell.Value = '0.9545'
$cell.Style = "Normal"

$cell = $ws.Range('E35')
$cell.NumberFormat = "@"
$cell.Value = '  +0.25%  '
$cell.Style = "Normal"

$cell = $ws.Range('D36')
$cell.NumberFormat = "@"
$cell.Value = '2.778'
$cell.Style = "Normal"

$cell = $ws.Range('E36')
$cell.NumberFormat = "@"
$cell.Value = '  -0.94%  '
$cell.Style = "Normal"

$cell = $ws.Range('E38')
$cell.NumberFormat = "@"
$cell.Value = '  -0.19%  '
$cell.Style = "Normal"

$cell = $ws.Range('D39')
$cell.NumberFormat = "@"
$cell.Value = '0.01588'
$cell.Style = "Normal"

$cell = $ws.Range('E39')
$cell.NumberFormat = "@"
$cell.Value = '  -0.61%  '
$cell.Style = "Normal"

$cell = $ws.Range('D40')
$cell.NumberFormat = "@"
$cell.Value = '5.882'
$cell.Style = "Normal"

$cell = $ws.Range('E40')
$cell.NumberFormat = "@"
$cell.Value = '  -1.12%  '
$cell.Style = "Normal"

$cell = $ws.Range('D42')
$cell.NumberFormat = "@"
$cell.Value = '0.8331'
$cell.Style = "Normal"

$cell = $ws.Range('E42')
$cell.NumberFormat = "@"
$cell.Value = '  -1.54%  '
$cell.Style = "Normal"

$cell = $ws.Range('D43')
$cell.NumberFormat = "@"
$cell.Value = '1.026.77'
$cell.Style = "Normal"

$cell = $ws.Range('E43')
$cell.NumberFormat = "@"
$cell.Value = '  -3.20%  '
$cell.Style = "Normal"

$cell = $ws.Range('D44')
$cell.NumberFormat = "@"
$cell.Value = '101.35'
$cell.Style = "Normal"

$cell = $ws.Range('E44')
$cell.NumberFormat = "@"
$cell.Value = '  -2.10%  '
$cell.Style = "Normal"

$cell = $ws.Range('D45')
$cell.NumberFormat = "@"
$cell.Value = '1.795.28'
$cell.Style = "Normal"

$cell = $ws.Range('E45')
$cell.NumberFormat = "@"
$cell.Value = '  +0.01%  '
$cell.Style = "Normal"

$cell = $ws.Range('D46')
$cell.NumberFormat = "@"
$cell.Value = '57.82'
$cell.Style = "Normal"

$cell = $ws.Range('E46')
$cell.NumberFormat = "@"
$cell.Value = '  -0.05%  '
$cell.Style = "Normal"

$cell = $ws.Range('D47')
$cell.NumberFormat = "@"
$cell.Value = '0.0₈105'
$cell.Style = "Normal"

$cell = $ws.Range('E47')
$cell.NumberFormat = "@"
$cell.Value = '  +1.62%  '
$cell.Style = "Normal"

$cell = $ws.Range('D48')
$cell.NumberFormat = "@"
$cell.Value = '1.002'
$cell.Style = "Normal"

$cell = $ws.Range('E48')
$cell.NumberFormat = "@"
$cell.Value = '  -0.64%  '
$cell.Style = "Normal"

$cell = $ws.Range('D49')
$cell.NumberFormat = "@"
$cell.Value = '8.031'
$cell.Style = "Normal"

$cell = $ws.Range('E49')
$cell.NumberFormat = "@"
$cell.Value = '  +0.48%  '
$cell.Style = "Normal"

$cell = $ws.Range('D50')
$cell.NumberFormat = "@"
$cell.Value = '0.4345'
$cell.Style = "Normal"

$cell = $ws.Range('E50')
$cell.NumberFormat = "@"
$cell.Value = '  -1.26%  '
$cell.Style = "Normal"

$cell = $ws.Range('D51')
$cell.NumberFormat = "@"
$cell.Value = '0.05199'
$cell.Style = "Normal"

$cell = $ws.Range('E51')
$cell.NumberFormat = "@"
$cell.Value = '  -3.98%  '
$cell.Style = "Normal"
